# Auto-generated edit script: apply scheduled market-data refresh values
# to the Gilgamesh Profits workbook (columns H-N per leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 310.16666
$ws.Range("I11").Value = 310.16666
$ws.Range("K11").Value = 310.16666
$ws.Range("M11").Value = -170.16666
$ws.Range("H28").Value = 1772.8823
$ws.Range("J28").Value = 5249.25
$ws.Range("L28").Value = 5249.25
$ws.Range("N28").Value = -6219.25
$ws.Range("H33").Value = 181.81818
$ws.Range("I33").Value = 191.1
$ws.Range("J33").Value = 89
$ws.Range("K33").Value = 191.1
$ws.Range("L33").Value = 89
$ws.Range("M33").Value = 37.90000000000001
$ws.Range("N33").Value = -547
$ws.Range("H43").Value = 2505.0833
$ws.Range("I43").Value = 3612
$ws.Range("J43").Value = 1951.625
$ws.Range("K43").Value = 3612
$ws.Range("L43").Value = 1951.625
$ws.Range("M43").Value = -3543
$ws.Range("N43").Value = -2089.625
$ws.Range("H103").Value = 2325
$ws.Range("I103").Value = 2528.625
$ws.Range("J103").Value = 1999.2
$ws.Range("K103").Value = 7585.875
$ws.Range("L103").Value = 5997.6
$ws.Range("M103").Value = -6999.875
$ws.Range("N103").Value = -7169.6
$ws.Range("H113").Value = 3767.2354
$ws.Range("I113").Value = 3161
$ws.Range("K113").Value = 3161
$ws.Range("M113").Value = 93
$ws.Range("H131").Value = 1002927.4
$ws.Range("I131").Value = 1431810.6
$ws.Range("K131").Value = 4295431.800000001
$ws.Range("M131").Value = -4290391.800000001
$ws.Range("H132").Value = 4388.227
$ws.Range("I132").Value = 4767.675
$ws.Range("J132").Value = 593.75
$ws.Range("K132").Value = 14303.025
$ws.Range("L132").Value = 1781.25
$ws.Range("M132").Value = -11773.025
$ws.Range("N132").Value = -6841.25
$ws.Range("H138").Value = 319815.88
$ws.Range("I138").Value = 3573.5
$ws.Range("J138").Value = 493439.12
$ws.Range("K138").Value = 10720.5
$ws.Range("L138").Value = 1480317.36
$ws.Range("M138").Value = -5580.5
$ws.Range("N138").Value = -1490597.36
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3813.8928
$ws.Range("I32").Value = 3877.8364
$ws.Range("K32").Value = 3877.8364
$ws.Range("M32").Value = -3590.8364
$ws.Range("H45").Value = 17968.258
$ws.Range("I45").Value = 22447.174
$ws.Range("J45").Value = 5091.375
$ws.Range("K45").Value = 22447.174
$ws.Range("L45").Value = 5091.375
$ws.Range("M45").Value = -22070.174
$ws.Range("N45").Value = -5845.375
$ws.Range("H102").Value = 6706.5
$ws.Range("J102").Value = 4757.7144
$ws.Range("L102").Value = 4757.7144
$ws.Range("N102").Value = -8001.7144
$ws.Range("H113").Value = 63548.5
$ws.Range("J113").Value = 63548.5
$ws.Range("L113").Value = 63548.5
$ws.Range("N113").Value = -72226.5
$ws.Range("H122").Value = 4139.852
$ws.Range("I122").Value = 3552.2727
$ws.Range("K122").Value = 10656.8181
$ws.Range("M122").Value = -8206.8181
$ws.Range("H132").Value = 2185.3555
$ws.Range("I132").Value = 1669.9143
$ws.Range("K132").Value = 5009.742899999999
$ws.Range("M132").Value = -2479.742899999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 414904
$ws.Range("J53").Value = 414904
$ws.Range("L53").Value = 414904
$ws.Range("N53").Value = -416052
$ws.Range("H86").Value = 3324.5
$ws.Range("I86").Value = 2870.9285
$ws.Range("K86").Value = 2870.9285
$ws.Range("M86").Value = -1747.9285
$ws.Range("H89").Value = 3324.5
$ws.Range("I89").Value = 2870.9285
$ws.Range("K89").Value = 14354.6425
$ws.Range("M89").Value = -8738.6425
$ws.Range("H99").Value = 4258.1875
$ws.Range("I99").Value = 3761
$ws.Range("K99").Value = 3761
$ws.Range("M99").Value = -2263
$ws.Range("H134").Value = 1569.9546
$ws.Range("I134").Value = 1182.75
$ws.Range("K134").Value = 3548.25
$ws.Range("M134").Value = -1013.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1742.4736
$ws.Range("I16").Value = 1684.7037
$ws.Range("J16").Value = 1884.2727
$ws.Range("K16").Value = 1684.7037
$ws.Range("L16").Value = 1884.2727
$ws.Range("M16").Value = -1397.7037
$ws.Range("N16").Value = -2458.2727
$ws.Range("H31").Value = 4748.913
$ws.Range("I31").Value = 4675
$ws.Range("K31").Value = 4675
$ws.Range("M31").Value = -4380
$ws.Range("H34").Value = 4748.913
$ws.Range("I34").Value = 4675
$ws.Range("K34").Value = 4675
$ws.Range("M34").Value = -4473
$ws.Range("H62").Value = 12517401
$ws.Range("J62").Value = 33664.668
$ws.Range("L62").Value = 33664.668
$ws.Range("N62").Value = -34912.668
$ws.Range("H65").Value = 12517401
$ws.Range("J65").Value = 33664.668
$ws.Range("L65").Value = 168323.34
$ws.Range("N65").Value = -174563.34
$ws.Range("H113").Value = 1742.4736
$ws.Range("I113").Value = 1684.7037
$ws.Range("J113").Value = 1884.2727
$ws.Range("K113").Value = 1684.7037
$ws.Range("L113").Value = 1884.2727
$ws.Range("M113").Value = 485.2963
$ws.Range("N113").Value = -6224.2727
$ws.Range("H134").Value = 4713.2593
$ws.Range("I134").Value = 4497.4546
$ws.Range("K134").Value = 13492.3638
$ws.Range("M134").Value = -10957.3638
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 3317.5334
$ws.Range("I21").Value = 2666.8333
$ws.Range("J21").Value = 3751.3333
$ws.Range("K21").Value = 8000.499899999999
$ws.Range("L21").Value = 11253.9999
$ws.Range("M21").Value = -7827.499899999999
$ws.Range("N21").Value = -11599.9999
$ws.Range("H26").Value = 1162.5
$ws.Range("I26").Value = 301
$ws.Range("J26").Value = 1285.5714
$ws.Range("K26").Value = 903
$ws.Range("L26").Value = 3856.7142
$ws.Range("M26").Value = -615
$ws.Range("N26").Value = -4432.7142
$ws.Range("H138").Value = 3562
$ws.Range("I138").Value = 2055.625
$ws.Range("J138").Value = 6574.75
$ws.Range("K138").Value = 6166.875
$ws.Range("L138").Value = 19724.25
$ws.Range("M138").Value = -1026.875
$ws.Range("N138").Value = -30004.25
$ws.Range("H139").Value = 2727.8
$ws.Range("I139").Value = 1850.5454
$ws.Range("K139").Value = 5551.6362
$ws.Range("M139").Value = -411.6361999999999
$ws.Range("H141").Value = 13376.8
$ws.Range("I141").Value = 13376.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 40130.39999999999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -34950.39999999999
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 928.6
$ws.Range("I22").Value = 513
$ws.Range("J22").Value = 1344.2
$ws.Range("K22").Value = 513
$ws.Range("L22").Value = 1344.2
$ws.Range("M22").Value = -218
$ws.Range("N22").Value = -1934.2
$ws.Range("H27").Value = 928.6
$ws.Range("I27").Value = 513
$ws.Range("J27").Value = 1344.2
$ws.Range("K27").Value = 513
$ws.Range("L27").Value = 1344.2
$ws.Range("M27").Value = -406
$ws.Range("N27").Value = -1558.2
$ws.Range("H40").Value = 42869.965
$ws.Range("I40").Value = 44134.383
$ws.Range("K40").Value = 44134.383
$ws.Range("M40").Value = -43998.383
$ws.Range("H46").Value = 2071
$ws.Range("I46").Value = 1844.5
$ws.Range("J46").Value = 2135.7144
$ws.Range("K46").Value = 1844.5
$ws.Range("L46").Value = 2135.7144
$ws.Range("M46").Value = -1656.5
$ws.Range("N46").Value = -2511.7144
$ws.Range("H55").Value = 493.5
$ws.Range("I55").Value = 429.8
$ws.Range("K55").Value = 429.8
$ws.Range("M55").Value = -256.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H113").Value = 649.0526
$ws.Range("I113").Value = 587.9231
$ws.Range("J113").Value = 781.5
$ws.Range("K113").Value = 1763.7693
$ws.Range("L113").Value = 2344.5
$ws.Range("M113").Value = 406.2307000000001
$ws.Range("N113").Value = -6684.5
$ws.Range("H126").Value = 2798
$ws.Range("I126").Value = 2785.625
$ws.Range("J126").Value = 2897
$ws.Range("K126").Value = 8356.875
$ws.Range("L126").Value = 8691
$ws.Range("M126").Value = -5886.875
$ws.Range("N126").Value = -13631
$ws.Range("H132").Value = 9262144
$ws.Range("I132").Value = 12823367
$ws.Range("J132").Value = 2966.9
$ws.Range("K132").Value = 38470101
$ws.Range("L132").Value = 8900.700000000001
$ws.Range("M132").Value = -38467571
$ws.Range("N132").Value = -13960.7
